# Update "江西-漫展信息" workbook:
#  - Bump the "想去人数" (interest count, column F) figures on sheets
#    "展览" (1) and "全部类型" (4) for most of the existing rows.
#  - Insert a brand-new convention (南昌·CM03动漫游戏博览会, 2024-08-17) as
#    row 40, pushing the two rows that used to be 40/41 down to 41/42
#    (their own F counts also tick up), growing the used range from
#    A1:I41 to A1:I42.
# Sheets "演出" (2) and "本地生活" (3) only contain a header row and are
# left untouched.

function Set-TextCell {
    # Force a literal-text write so values that look like dates (e.g.
    # "2024-08-17") are not silently reinterpreted as date serials.
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# F-column (想去人数) row -> (old, new) bumps, identical on both sheets.
$fUpdates = @{
    2  = 259
    3  = 1341
    4  = 148
    6  = 226
    7  = 96
    9  = 181
    10 = 129
    11 = 4521
    12 = 6788
    13 = 41
    14 = 58
    16 = 568
    18 = 4122
    19 = 486
    21 = 55
    22 = 2696
    24 = 547
    25 = 165
    26 = 352
    27 = 358
    29 = 223
    30 = 34
    31 = 1623
    32 = 1018
    33 = 62
    34 = 131
    35 = 79
    36 = 542
    37 = 496
    39 = 88
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # 1) Bump the F-column interest counts for the untouched rows.
    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # 2) Make room for the new row 40 by rotating the current rows
    #    40 and 41 down to 41 and 42 (content only; the A-column serial
    #    numbers 39/40/41 stay put since they are positional, not tied
    #    to the row's event).

    # 2a) Give row 42 row 41's style + value via copy, fix value after.
    $ws.Cells.Item(41, 1).Copy($ws.Cells.Item(42, 1))
    $ws.Cells.Item(42, 1).Value = 41

    # 2b) Old row 41 (哥布林展) content -> row 42.
    $ws.Range("B41:I41").Copy($ws.Range("B42:I42"))

    # 2c) Old row 40 (龙年动漫展) content -> row 41.
    $ws.Range("B40:I40").Copy($ws.Range("B41:I41"))

    # 3) Fix up the F values on the two rows that were shifted down.
    $ws.Cells.Item(41, 6).Value = 639
    $ws.Cells.Item(42, 6).Value = 10

    # 4) Overwrite row 40 with the brand-new convention entry.
    Set-TextCell $ws.Cells.Item(40, 2) "2024-08-17"
    $ws.Cells.Item(40, 3).Value = "南昌·CM03动漫游戏博览会"
    $ws.Cells.Item(40, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(40, 5).Value = "2024.08.17 09:00-08.18 17:00"
    $ws.Cells.Item(40, 6).Value = 55
    $ws.Cells.Item(40, 7).Value = 55
    $ws.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89295"
    $ws.Cells.Item(40, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/LoSq5kzH1719840007767.png"
}
